$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, centered, bordered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate I0/IF data for rows 2-70
$ijData = @{
    2 = @(7, 7)
    3 = @(8, 8)
    4 = @(7, 8)
    5 = @(6, 7)
    6 = @(6, 7)
    7 = @(9, 9)
    8 = @(7, 7)
    9 = @(8, 8)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(5, 5)
    13 = @(6, 6)
    14 = @(6, 6)
    15 = @(7, 7)
    16 = @(9, 9)
    17 = @(8, 8)
    18 = @(9, 9)
    19 = @(8, 8)
    20 = @(7, 8)
    21 = @(9, 9)
    22 = @(7, 7)
    23 = @(8, 8)
    24 = @(7, 8)
    25 = @(8, 8)
    26 = @(7, 8)
    27 = @(6, 7)
    28 = @(5, 6)
    29 = @(7, 7)
    30 = @(8, 8)
    31 = @(7, 8)
    32 = @(7, 7)
    33 = @(7, 7)
    34 = @(8, 8)
    35 = @(5, 6)
    36 = @(8, 8)
    37 = @(7, 8)
    38 = @(8, 8)
    39 = @(5, 5)
    40 = @(6, 6)
    41 = @(8, 8)
    42 = @(8, 8)
    43 = @(6, 6)
    44 = @(7, 7)
    45 = @(7, 8)
    46 = @(9, 9)
    47 = @(8, 8)
    48 = @(6, 6)
    49 = @(8, 8)
    50 = @(8, 8)
    51 = @(8, 8)
    52 = @(7, 8)
    53 = @(10, 10)
    54 = @(9, 9)
    55 = @(8, 8)
    56 = @(9, 9)
    57 = @(7, 7)
    58 = @(8, 8)
    59 = @(5, 6)
    60 = @(8, 8)
    61 = @(8, 8)
    62 = @(6, 7)
    63 = @(8, 8)
    64 = @(8, 8)
    65 = @(7, 7)
    66 = @(4, 4)
    67 = @(5, 5)
    68 = @(4, 4)
    69 = @(6, 6)
    70 = @(6, 6)
}

foreach ($r in $ijData.Keys) {
    $pair = $ijData[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
